# Estadisticos Segundo Parcial 26 Mayo
#
# 1) "Calificaciones": fill in the previously-missing 2P grade (column J) for
#    subject "Conciencia historica 1..." for every student, and update the
#    Final grade (column X) for that same subject for the students whose
#    final recalculated once the 2P grade was entered.
# 2) "Totales": update the aggregate row for subject "Conciencia historica 1..."
#    (row 6) to reflect that all 31 students now passed.
# 3) "Rescatables": add four new rows for VIVANCO VIVANCO LUIS AARON (one per
#    subject he needs to retake) ahead of the pre-existing three rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Calificaciones sheet - columns J (2P) and X (Final) for subject "Conciencia
#    historica 1: Perspectivas del Mexico antiguo en los contextos globales"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Calificaciones")

$jValues = @{
    4  = 10; 5  = 10; 6  = 7;  7  = 10; 8  = 9;  9  = 7;  10 = 7;
    11 = 7;  12 = 10; 13 = 10; 14 = 10; 15 = 10; 16 = 6;  17 = 5;
    18 = 10; 19 = 8;  20 = 10; 21 = 10; 22 = 10; 23 = 10; 24 = 10;
    25 = 10; 26 = 10; 27 = 10; 28 = 10; 29 = 10; 30 = 10; 31 = 10;
    32 = 9;  33 = 7;  34 = 7
}

$xValues = @{
    4 = 10; 7 = 10; 10 = 9; 11 = 9; 12 = 9; 13 = 9; 14 = 9; 16 = 8;
    17 = 6; 21 = 9; 22 = 9; 23 = 10; 25 = 8; 26 = 10; 27 = 10; 31 = 8;
    33 = 6
}

foreach ($row in $jValues.Keys) {
    $ws1.Range("J$row").Value = $jValues[$row]
}

foreach ($row in $xValues.Keys) {
    $ws1.Range("X$row").Value = $xValues[$row]
}

# ---------------------------------------------------------------------------
# 2) Totales sheet - row 6 is the summary for the same subject
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Totales")
$ws3.Range("D6").Value = 31
$ws3.Range("E6").Value = 0
$ws3.Range("F6").Value = 100
$ws3.Range("G6").Value = 0
$ws3.Range("H6").Value = 9

# ---------------------------------------------------------------------------
# 3) Rescatables sheet - insert 4 new rows (VIVANCO VIVANCO LUIS AARON) before
#    the existing 3 rows, pushing them down to rows 6-8. Rewritten directly
#    (rather than via Rows.Insert) so no stray styles are introduced.
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Rescatables")

# Preserve + shift the 3 pre-existing rows (previously rows 2-4) down to 6-8
$ws5.Range("A8").Value = 23330051920311
$ws5.Range("B8").Value = "TEMOXTLE"
$ws5.Range("C8").Value = "GARCIA"
$ws5.Range("D8").Value = "HUGO ANTONIO"
$ws5.Range("E8").Value = "Reacciones químicas: conservación de la materia en la formación de nuevas substancias."
$ws5.Range("F8").Value = "Nativitas Sandoval Liliana Soledad"
$ws5.Range("G8").Value = 5

$ws5.Range("A7").Value = 23330051920164
$ws5.Range("B7").Value = "LARRACILLA"
$ws5.Range("C7").Value = "GOMEZ"
$ws5.Range("D7").Value = "MONICA"
$ws5.Range("E7").Value = "IMPLEMENTA BASE DE DATOS RELACIONALES EN UN SISTEMA DE INFORMACIÓN"
$ws5.Range("F7").Value = "Rodriguez Roman Marisol"
$ws5.Range("G7").Value = 5

$ws5.Range("A6").Value = 23330051920155
$ws5.Range("B6").Value = "CRUZ"
$ws5.Range("C6").Value = "NIEVES"
$ws5.Range("D6").Value = "ESTRELLA ESMERALDA"
$ws5.Range("E6").Value = "Ingles IV"
$ws5.Range("F6").Value = "Avila Coronado Julieta"
$ws5.Range("G6").Value = 5

# New student: VIVANCO VIVANCO LUIS AARON, one row per subject to retake
$ws5.Range("A2").Value = 23330051920313
$ws5.Range("B2").Value = "VIVANCO"
$ws5.Range("C2").Value = "VIVANCO"
$ws5.Range("D2").Value = "LUIS AARON"
$ws5.Range("E2").Value = "IMPLEMENTA BASE DE DATOS RELACIONALES EN UN SISTEMA DE INFORMACIÓN"
$ws5.Range("F2").Value = "Rodriguez Roman Marisol"
$ws5.Range("G2").Value = 5

$ws5.Range("A3").Value = 23330051920313
$ws5.Range("B3").Value = "VIVANCO"
$ws5.Range("C3").Value = "VIVANCO"
$ws5.Range("D3").Value = "LUIS AARON"
$ws5.Range("E3").Value = "Ingles IV"
$ws5.Range("F3").Value = "Avila Coronado Julieta"
$ws5.Range("G3").Value = 5

$ws5.Range("A4").Value = 23330051920313
$ws5.Range("B4").Value = "VIVANCO"
$ws5.Range("C4").Value = "VIVANCO"
$ws5.Range("D4").Value = "LUIS AARON"
$ws5.Range("E4").Value = "Reacciones químicas: conservación de la materia en la formación de nuevas substancias."
$ws5.Range("F4").Value = "Nativitas Sandoval Liliana Soledad"
$ws5.Range("G4").Value = 5

$ws5.Range("A5").Value = 23330051920313
$ws5.Range("B5").Value = "VIVANCO"
$ws5.Range("C5").Value = "VIVANCO"
$ws5.Range("D5").Value = "LUIS AARON"
$ws5.Range("E5").Value = "Temas selectos de matemáticas I"
$ws5.Range("F5").Value = "Ortega Medina Angel Gaspar"
$ws5.Range("G5").Value = 5
